$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (often containing multiple dots, or
# subscript digits) that must stay text. Force text formatting on
# each target D-cell first so Excel does not auto-convert the new
# value into a floating point number (which would both change the
# stored type and introduce binary floating point rounding noise).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.176.64'
$ws.Range("E2").Value = '  +6.15%  '
$ws.Range("D3").Value = '3.108.13'
$ws.Range("E3").Value = '  +4.18%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '584.61'
$ws.Range("E5").Value = '  +3.16%  '
$ws.Range("D6").Value = '144.72'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.102.47'
$ws.Range("E8").Value = '  +4.19%  '
$ws.Range("E9").Value = '  +1.98%  '
$ws.Range("E10").Value = '  +13.35%  '
$ws.Range("D11").Value = '5.77'
$ws.Range("E11").Value = '  +7.12%  '
$ws.Range("D12").Value = '0.467'
$ws.Range("E12").Value = '  +3.65%  '
$ws.Range("E13").Value = '  +7.34%  '
$ws.Range("E14").Value = '  +5.21%  '
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("D16").Value = '3.622.97'
$ws.Range("E16").Value = '  +4.22%  '
$ws.Range("D17").Value = '7.15'
$ws.Range("E17").Value = '  +1.13%  '
$ws.Range("D18").Value = '63.099.90'
$ws.Range("E18").Value = '  +6.04%  '
$ws.Range("D19").Value = '3.109.16'
$ws.Range("E19").Value = '  +4.25%  '
$ws.Range("D20").Value = '466.48'
$ws.Range("E20").Value = '  +6.56%  '
$ws.Range("E21").Value = '  +3.62%  '
$ws.Range("D22").Value = '0.726'
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("E23").Value = '  +6.85%  '
$ws.Range("D24").Value = '13.31'
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").Value = '81.99'
$ws.Range("E25").Value = '  +2.44%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  +9.75%  '
$ws.Range("E28").Value = '  -0.75%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  +4.65%  '
$ws.Range("D31").Value = '6.87'
$ws.Range("E31").Value = '  +10.12%  '
$ws.Range("D32").Value = '26.92'
$ws.Range("E32").Value = '  +4.38%  '
$ws.Range("E33").Value = '  +3.41%  '
$ws.Range("D34").Value = '0.0₃0861'
$ws.Range("E34").Value = '  +10.60%  '
$ws.Range("E35").Value = '  +15.71%  '
$ws.Range("E36").Value = '  +6.96%  '
$ws.Range("D37").Value = '6.07'
$ws.Range("E37").Value = '  +2.55%  '
$ws.Range("E38").Value = '  +19.24%  '
$ws.Range("D39").Value = '50.57'
$ws.Range("E39").Value = '  +3.89%  '
$ws.Range("D40").Value = '436.78'
$ws.Range("E40").Value = '  +9.36%  '
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("D42").Value = '2.916.93'
$ws.Range("E42").Value = '  +6.65%  '
$ws.Range("E43").Value = '  +4.74%  '
$ws.Range("E44").Value = '  +10.92%  '
$ws.Range("D46").Value = '2.16'
$ws.Range("E46").Value = '  +7.29%  '
$ws.Range("D48").Value = '34.79'
$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("D49").Value = '122.92'
$ws.Range("E49").Value = '  +0.57%  '
$ws.Range("E50").Value = '  +0.72%  '
$ws.Range("D51").Value = '24.59'
$ws.Range("E51").Value = '  +5.32%  '

# Restore the cells default (unstyled) formatting so the saved
# workbook styling matches the original file (no leftover @ text
# format hanging around on these cells).
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
